# Update the Engineering_data values on Sheet1 (row 4: box_hole_depth)
# to fit the adjusted data model, then move the active selection to C4
# (matching the cursor position the workbook was saved with).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$ws.Range("B4").Value = 30
$ws.Range("C4").Value = 28.5
$ws.Range("D4").Value = 31.5

$ws.Range("C4").Select()
